$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.824.93"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").Value = "3.924.28"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.44"
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.68"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -1.00%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.738"
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("E10").Value = "  +3.39%  "

$ws.Range("E11").Value = "  -3.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.24"
$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.46"
$ws.Range("E13").Value = "  -1.41%  "

$ws.Range("D14").Value = "4.563.32"

$ws.Range("D15").Value = "3.920.11"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.32"
$ws.Range("E16").Value = "  -3.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.98"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("E19").Value = "  +2.41%  "

$ws.Range("D20").Value = "68.940.51"
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.76"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.49"
$ws.Range("E22").Value = "  +4.05%  "

$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.45"
$ws.Range("E24").Value = "  +14.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "89.38"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("E26").Value = "  +3.12%  "

$ws.Range("E27").Value = "  -3.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.11"
$ws.Range("E28").Value = "  -4.54%  "

$ws.Range("E29").Value = "  -3.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "710.61"
$ws.Range("E30").Value = "  +3.58%  "

$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.490"
$ws.Range("E34").Value = "  +33.36%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.19"
$ws.Range("E35").Value = "  +9.40%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0887"
$ws.Range("E36").Value = "  -4.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "61.47"
$ws.Range("E37").Value = "  +3.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "40.95"
$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.97"
$ws.Range("E42").Value = "  +3.34%  "

$ws.Range("E43").Value = "  +1.88%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("D46").Value = "0.0₆0376"
$ws.Range("E46").Value = "  +14.73%  "

$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("E48").Value = "  +8.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.02"
$ws.Range("E49").Value = "  +6.21%  "

$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("E51").Value = "  -2.60%  "
